$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell's value as literal text, preserving the original
# (unstyled / General) number format even for values that look numeric
# (e.g. '297.03', '1.00') so Excel doesn't silently coerce them to numbers.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Row 2
Set-TextValue $ws.Range("D2") '42.848.34'
Set-TextValue $ws.Range("E2") '  -7.58%  '
# Row 3
Set-TextValue $ws.Range("D3") '2.525.54'
Set-TextValue $ws.Range("E3") '  -3.39%  '
# Row 4
Set-TextValue $ws.Range("D4") '1.00'
Set-TextValue $ws.Range("E4") '  -0.04%  '
# Row 5
Set-TextValue $ws.Range("D5") '297.03'
Set-TextValue $ws.Range("E5") '  -3.56%  '
# Row 6
Set-TextValue $ws.Range("D6") '93.28'
Set-TextValue $ws.Range("E6") '  -6.81%  '
# Row 7
Set-TextValue $ws.Range("E7") '  -5.21%  '
# Row 8
Set-TextValue $ws.Range("D8") '1.00'
Set-TextValue $ws.Range("E8") '  +0.07%  '
# Row 9
Set-TextValue $ws.Range("D9") '0.547'
Set-TextValue $ws.Range("E9") '  -5.89%  '
# Row 10
Set-TextValue $ws.Range("D10") '36.24'
Set-TextValue $ws.Range("E10") '  -7.89%  '
# Row 11
Set-TextValue $ws.Range("D11") '0.0802'
Set-TextValue $ws.Range("E11") '  -5.02%  '
# Row 12
Set-TextValue $ws.Range("D12") '7.57'
Set-TextValue $ws.Range("E12") '  -7.00%  '
# Row 13
Set-TextValue $ws.Range("E13") '  +0.46%  '
# Row 14
Set-TextValue $ws.Range("D14") '2.910.59'
Set-TextValue $ws.Range("E14") '  -3.49%  '
# Row 15
Set-TextValue $ws.Range("D15") '2.522.99'
Set-TextValue $ws.Range("E15") '  -3.60%  '
# Row 16
Set-TextValue $ws.Range("E16") '  -6.56%  '
# Row 17
Set-TextValue $ws.Range("D17") '14.07'
Set-TextValue $ws.Range("E17") '  -6.30%  '
# Row 18
Set-TextValue $ws.Range("D18") '42.900.05'
Set-TextValue $ws.Range("E18") '  -7.79%  '
# Row 19
Set-TextValue $ws.Range("B19") 'Uniswap'
Set-TextValue $ws.Range("C19") 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue $ws.Range("D19") '6.55'
Set-TextValue $ws.Range("E19") '  -3.31%  '
# Row 20
Set-TextValue $ws.Range("B20") 'ShibaInu'
Set-TextValue $ws.Range("C20") 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue $ws.Range("D20") '0.0₃0962'
Set-TextValue $ws.Range("E20") '  -5.10%  '
# Row 21
Set-TextValue $ws.Range("E21") '  -5.71%  '
# Row 22
Set-TextValue $ws.Range("D22") '72.81'
Set-TextValue $ws.Range("E22") '  +1.22%  '
# Row 23
Set-TextValue $ws.Range("D23") '258.48'
Set-TextValue $ws.Range("E23") '  -6.21%  '
# Row 24
Set-TextValue $ws.Range("D24") '2.90'
Set-TextValue $ws.Range("E24") '  -4.41%  '
# Row 25
Set-TextValue $ws.Range("D25") '2.16'
Set-TextValue $ws.Range("E25") '  -2.29%  '
# Row 26
Set-TextValue $ws.Range("D26") '28.91'
Set-TextValue $ws.Range("E26") '  -2.29%  '
# Row 27
Set-TextValue $ws.Range("E27") '  +0.15%  '
# Row 28
Set-TextValue $ws.Range("D28") '9.95'
Set-TextValue $ws.Range("E28") '  -6.58%  '
# Row 29
Set-TextValue $ws.Range("B29") 'InjectiveProtocol'
Set-TextValue $ws.Range("C29") 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws.Range("D29") '36.77'
Set-TextValue $ws.Range("E29") '  -5.03%  '
# Row 30
Set-TextValue $ws.Range("B30") 'Toncoin'
Set-TextValue $ws.Range("C30") 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue $ws.Range("D30") '2.13'
Set-TextValue $ws.Range("E30") '  -4.44%  '
# Row 31
Set-TextValue $ws.Range("D31") '5.95'
Set-TextValue $ws.Range("E31") '  -7.42%  '
# Row 32
Set-TextValue $ws.Range("E32") '  -4.53%  '
# Row 33
Set-TextValue $ws.Range("D33") '2.20'
Set-TextValue $ws.Range("E33") '  -2.20%  '
# Row 34
Set-TextValue $ws.Range("D34") '151.39'
Set-TextValue $ws.Range("E34") '  -0.60%  '
# Row 35
Set-TextValue $ws.Range("D35") '2.76'
Set-TextValue $ws.Range("E35") '  -3.03%  '
# Row 36
Set-TextValue $ws.Range("E36") '  -4.88%  '
# Row 37
Set-TextValue $ws.Range("E37") '  -6.03%  '
# Row 38
Set-TextValue $ws.Range("E38") '  -3.93%  '
# Row 39
Set-TextValue $ws.Range("D39") '23.71'
Set-TextValue $ws.Range("E39") '  -0.71%  '
# Row 40
Set-TextValue $ws.Range("D40") '16.26'
Set-TextValue $ws.Range("E40") '  +1.71%  '
# Row 41
Set-TextValue $ws.Range("D41") '3.43'
Set-TextValue $ws.Range("E41") '  -5.38%  '
# Row 42
Set-TextValue $ws.Range("D42") '0.0308'
Set-TextValue $ws.Range("E42") '  -7.02%  '
# Row 43
Set-TextValue $ws.Range("E43") '  -5.74%  '
# Row 44
Set-TextValue $ws.Range("D44") '2.016.35'
Set-TextValue $ws.Range("E44") '  -5.65%  '
# Row 45
Set-TextValue $ws.Range("D45") '1.00'
Set-TextValue $ws.Range("E45") '  +0.16%  '
# Row 46
Set-TextValue $ws.Range("D46") '85.70'
Set-TextValue $ws.Range("E46") '  -9.37%  '
# Row 48
Set-TextValue $ws.Range("D48") '8.87'
Set-TextValue $ws.Range("E48") '  -6.93%  '
# Row 49
Set-TextValue $ws.Range("D49") '2.765.91'
Set-TextValue $ws.Range("E49") '  -3.70%  '
# Row 50
Set-TextValue $ws.Range("D50") '102.80'
Set-TextValue $ws.Range("E50") '  -6.36%  '
# Row 51
Set-TextValue $ws.Range("B51") 'Algorand'
Set-TextValue $ws.Range("C51") 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws.Range("D51") '0.187'
Set-TextValue $ws.Range("E51") '  -7.48%  '
